$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values that look like plain numbers need an explicit text
# NumberFormat first so Excel keeps them as strings (matches source data,
# which stores every Price cell as text, e.g. "22.25", "1.785.83").

$ws.Range("D2").Value = "27.136.18"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "1.565.51"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.63"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.487"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.29"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.787.84"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").Value = "1.562.54"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "27.148.22"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.68"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("D19").Value = "0.0₃0685"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.25"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.35"
$ws.Range("E23").Value = "  -4.01%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.06"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  -7.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.94"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").Value = "1.392.49"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.942"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.517"
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.993"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.31"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "1.701.10"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.57"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").Value = "0.0₇0981"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0493"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0949"
$ws.Range("E51").Value = "  -1.81%  "
